$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new weekly scoreboard entry in row 72 (just below the last existing
# data row, 71). This mirrors a manual row of data being typed into the sheet.
$row = 72

$ws.Cells.Item($row, 1).Value = "Stev en"          # Participant (A)

# Copy the date format from the cell above so the new date cell reuses the
# existing style (rather than Excel registering a brand-new number format).
$ws.Range("B" + ($row - 1)).Copy()
$ws.Range("B" + $row).PasteSpecial(-4122)          # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 2).Value = 45462               # Date (B) -> 6/19/2024

$ws.Cells.Item($row, 3).Value = "Workout"          # Workout Type (C)
$ws.Cells.Item($row, 4).Value = 49                 # Total Duration (D)
$ws.Cells.Item($row, 5).Value = 0                  # Total Distance (E)
$ws.Cells.Item($row, 6).Value = 0                  # Total Elevation (F)
$ws.Cells.Item($row, 7).Value = 1                  # Zone 1 (G)
$ws.Cells.Item($row, 8).Value = 17                 # Zone 2 (H)
$ws.Cells.Item($row, 9).Value = 25                 # Zone 3 (I)
$ws.Cells.Item($row, 10).Value = 7                 # Zone 4 (J)
$ws.Cells.Item($row, 11).Value = 0                 # Zone 5 (K)
$ws.Cells.Item($row, 12).Value = "Agile Antelope"  # Workout Level (L)
$ws.Cells.Item($row, 13).Value = 2                 # Week (M)

$ws.Range("M72").Select()

$wb.Save()
